$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Workflow State" header to "Review State"
$ws.Range("I1").Value = "Review State"

# Update selection to match the recorded state after the edit
$ws.Range("I2").Select()
